$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
  @(1.02, 1.025748667024775, 1.029894664102174, 1.049995176591597, 1.054727962127409, 1.032777525602488, 1.030915954402657, 1.032707250685951, 1.052750685033849, 1.05747037417289, 1.014305742470162),
  @(1.02, 1.026875634722425, 1.030712288360634, 1.051169488970474, 1.055953782723888, 1.033027235030032, 1.03168147338287, 1.033333090759105, 1.053736421677713, 1.058508441888384, 1.014561686930693),
  @(1.02, 1.027604309608448, 1.031240596846466, 1.051929834866339, 1.056747170747412, 1.03318681678612, 1.032175756601925, 1.03373665747303, 1.054374171086916, 1.059179775932418, 1.01472686667787),
  @(1.02, 1.027910514630509, 1.031462518553196, 1.052249602051823, 1.057080760227404, 1.033253426934924, 1.032383300433412, 1.033905983574127, 1.054642260831451, 1.059461918373407, 1.014796204629282),
  @(1.02, 1.0279619202251, 1.031499769681203, 1.052303299307959, 1.057136774287512, 1.03326458304193, 1.03241813314791, 1.033934394632577, 1.05468727312879, 1.059509286294703, 1.014807840704758),
  @(1.02, 1.027608401643584, 1.031243562879069, 1.051934107146966, 1.056751627991679, 1.033187708712441, 1.032178530805453, 1.03373892132615, 1.054377753392982, 1.059183546271701, 1.014727793581923),
  @(1.02, 1.026129645483037, 1.030171139245651, 1.05039194074173, 1.055142193249132, 1.032862329585144, 1.031174884731215, 1.032919045019715, 1.053083838187737, 1.05782127025483, 1.01439232993572),
  @(1.02, 1.023519622051221, 1.028275646822027, 1.0476781277619, 1.052307640077705, 1.032273672873009, 1.02939819418326, 1.031463630053278, 1.050803074751453, 1.055417922092529, 1.013797873507461),
  @(1.02, 1.02177664449283, 1.02700810791891, 1.045871320727614, 1.050418848943659, 1.031870943924535, 1.028208218574993, 1.030486144125803, 1.049282017056054, 1.053813709511066, 1.013399321658474),
  @(1.02, 1.021021191636198, 1.026458325535395, 1.045089505774476, 1.049601178475678, 1.031694114868636, 1.0276916271992, 1.030061165873886, 1.048623236890628, 1.053118581870441, 1.013226208445077),
  @(1.02, 1.020740470828007, 1.026253971700521, 1.044799185105279, 1.049297485828632, 1.031628065118462, 1.027499542202599, 1.029903050917875, 1.048378512340255, 1.052860304946288, 1.013161825436644),
  @(1.02, 1.020800691493811, 1.026297812643427, 1.044861456269848, 1.049362627758211, 1.031642249651567, 1.027540754151599, 1.029936978848408, 1.048431007707721, 1.052915709667018, 1.013175639482861),
  @(1.02, 1.020997989462572, 1.026441436446227, 1.045065506150182, 1.049576074614873, 1.031688662675479, 1.027675753478543, 1.030048101331834, 1.048603008386099, 1.053097234169579, 1.013220888180605),
  @(1.02, 1.021119536354386, 1.026529909220954, 1.045191238597897, 1.0497075897362, 1.031717210557219, 1.027758904492012, 1.030116533207276, 1.048708980423312, 1.053209067348229, 1.013248756647928),
  @(1.02, 1.021826765816939, 1.027044575535017, 1.045923218546893, 1.050473118862571, 1.031882627955581, 1.028242475058658, 1.030514312231528, 1.049325734800516, 1.053859832328346, 1.013410799261199),
  @(1.02, 1.022270194773246, 1.027367162863909, 1.046382514874395, 1.050953364303269, 1.031985735200354, 1.028545450881064, 1.030763367456859, 1.049712566802319, 1.054267906960009, 1.013512300175492),
  @(1.02, 1.022528768881379, 1.027555232929695, 1.04665046728594, 1.051233501655056, 1.032045640045894, 1.028722043866996, 1.030908471206584, 1.049938184779075, 1.054505882521568, 1.013571452086743),
  @(1.02, 1.022616924072618, 1.027619344768335, 1.046741841109368, 1.051329024411866, 1.032066026051222, 1.028782235885144, 1.030957919706844, 1.050015112247455, 1.054587018032483, 1.013591612572846),
  @(1.02, 1.022222626349036, 1.027332561574419, 1.046333231292115, 1.05090183663086, 1.031974697167321, 1.028512957656965, 1.030736663357597, 1.049671064908016, 1.05422412933432, 1.013501415453518),
  @(1.02, 1.020939893235589, 1.02639914668977, 1.045005416344612, 1.049513219116253, 1.031675005345541, 1.027636005069678, 1.030015385688058, 1.048552359154469, 1.053043781826073, 1.013207565806641),
  @(1.02, 1.020132739134134, 1.025811460415142, 1.0441710296988, 1.048640293030103, 1.031484450444352, 1.027083472757157, 1.029560390327601, 1.047848843357448, 1.052301212539918, 1.013022341781467),
  @(1.02, 1.020560688970983, 1.026123081101879, 1.044613310458633, 1.049103033764302, 1.031585668886206, 1.027376490646471, 1.029801734325769, 1.048221804237179, 1.052694904584637, 1.013120577100412),
  @(1.02, 1.022244120689646, 1.027348196684535, 1.046355500258472, 1.050925119710948, 1.031979685507957, 1.028527640340455, 1.030748730303073, 1.049689817872273, 1.054243910705193, 1.013506333950966),
  @(1.02, 1.024194890830733, 1.028766359319202, 1.048379284763158, 1.053040273434327, 1.032427668141756, 1.029858480130137, 1.031841158652659, 1.051392798708372, 1.056039589226758, 1.013951950077347)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $rowNum = 2 + $r
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $addr = $cols[$c] + $rowNum
        $ws.Range($addr).Value = $rowVals[$c]
    }
}

Write-Host "Done updating vm_pu values"
